$d = $word.ActiveDocument

$replacements = @(
    @("50×23=", "12×72="),
    @("65×12=", "97×22="),
    @("21×25=", "63×23="),
    @("15×56=", "31×39="),
    @("93×54=", "33×18="),
    @("44×72=", "37×72="),
    @("70×55=", "43×35="),
    @("79×28=", "11×39="),
    @("96×18=", "71×68="),
    @("71×99=", "46×61="),
    @("92×19=", "42×49="),
    @("45×65=", "36×71="),
    @("82×77=", "39×51="),
    @("14×31=", "70×21="),
    @("98×91=", "56×64="),
    @("58×93=", "70×96="),
    @("42×50=", "45×59="),
    @("25×17=", "55×71="),
    @("23×87=", "97×45="),
    @("98×60=", "19×95="),
    @("25×18=", "36×49="),
    @("33×44=", "99×33="),
    @("57×49=", "31×15="),
    @("63×13=", "37×74="),
    @("18×81=", "86×87=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
